$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT, even when it looks like a number
# (e.g. "250.86"), so Excel doesn't silently coerce it to a numeric
# cell. Forces the "@" text format just for the assignment, then
# restores the cell's style to Normal/General so no visible formatting
# change is left behind.
function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "41.480.42"
$ws.Range("E2").Value = "  -1.22%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.191.73"
$ws.Range("E3").Value = "  -1.58%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.15%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "250.86"
$ws.Range("E5").Value = "  -0.62%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -2.53%  "

# Row 7 - Solana
Set-TextValue $ws.Range("D7") "67.66"
$ws.Range("E7").Value = "  -5.94%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.10%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.617"
$ws.Range("E9").Value = "  +2.17%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "40.31"
$ws.Range("E10").Value = "  -0.40%  "

# Row 11 - OKB
Set-TextValue $ws.Range("D11") "59.68"
$ws.Range("E11").Value = "  +2.30%  "

# Row 12 - Dogecoin
$ws.Range("E12").Value = "  -3.17%  "

# Row 13 - Polkadot
Set-TextValue $ws.Range("D13") "7.05"
$ws.Range("E13").Value = "  -3.40%  "

# Row 14 - TRON
$ws.Range("E14").Value = "  -0.37%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "2.518.90"
$ws.Range("E15").Value = "  -1.58%  "

# Row 16 - Chainlink
Set-TextValue $ws.Range("D16") "14.51"
$ws.Range("E16").Value = "  -3.92%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.858"
$ws.Range("E17").Value = "  -3.77%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "2.177.96"
$ws.Range("E18").Value = "  -2.09%  "

# Row 19 - WrappedBTC
Set-TextValue $ws.Range("D19") "41.422.44"
$ws.Range("E19").Value = "  -1.23%  "

# Row 20 - ShibaInu
Set-TextValue $ws.Range("D20") "0.0₃0952"
$ws.Range("E20").Value = "  -1.90%  "

# Row 21 - was Litecoin, now Uniswap
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextValue $ws.Range("D21") "6.12"
$ws.Range("E21").Value = "  -2.95%  "

# Row 22 - was Uniswap, now Litecoin
$ws.Range("B22").Value = "Litecoin"
$ws.Range("C22").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws.Range("D22") "72.02"
$ws.Range("E22").Value = "  -1.25%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "231.33"
$ws.Range("E23").Value = "  -1.95%  "

# Row 24 - ImmutableX
Set-TextValue $ws.Range("D24") "2.07"
$ws.Range("E24").Value = "  -1.46%  "

# Row 25 - WEMIXToken
Set-TextValue $ws.Range("D25") "3.86"
$ws.Range("E25").Value = "  -6.04%  "

# Row 26 - was Dai, now Cosmos
$ws.Range("B26").Value = "Cosmos"
$ws.Range("C26").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-TextValue $ws.Range("D26") "11.40"
$ws.Range("E26").Value = "  -5.33%  "

# Row 27 - was Cosmos, now Dai
$ws.Range("B27").Value = "Dai"
$ws.Range("C27").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
Set-TextValue $ws.Range("D27") "1.00"
$ws.Range("E27").Value = "  +0.06%  "

# Row 28 - PancakeSwap
$ws.Range("E28").Value = "  -4.87%  "

# Row 29 - LEO
$ws.Range("E29").Value = "  -3.05%  "

# Row 30 - Toncoin
Set-TextValue $ws.Range("D30") "2.15"
$ws.Range("E30").Value = "  -2.14%  "

# Row 31 - Monero
Set-TextValue $ws.Range("D31") "167.14"
$ws.Range("E31").Value = "  -2.25%  "

# Row 32 - EthereumClassic
Set-TextValue $ws.Range("D32") "20.33"
$ws.Range("E32").Value = "  -2.84%  "

# Row 33 - Kaspa
Set-TextValue $ws.Range("D33") "0.121"
$ws.Range("E33").Value = "  -2.32%  "

# Row 34 - Hedera
Set-TextValue $ws.Range("D34") "0.0784"
$ws.Range("E34").Value = "  +5.06%  "

# Row 35 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D35") "5.79"
$ws.Range("E35").Value = "  +3.38%  "

# Row 36 - Stellar
$ws.Range("E36").Value = "  -2.18%  "

# Row 37 - RenderToken
$ws.Range("E37").Value = "  +3.69%  "

# Row 38 - InjectiveProtocol
Set-TextValue $ws.Range("D38") "26.06"
$ws.Range("E38").Value = "  -1.03%  "

# Row 39 - Filecoin
$ws.Range("E39").Value = "  -3.19%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.70%  "

# Row 41 - LidoDAOToken
$ws.Range("E41").Value = "  -3.02%  "

# Row 42 - FTXToken
Set-TextValue $ws.Range("D42") "5.21"
$ws.Range("E42").Value = "  +7.06%  "

# Row 43 - THORChain
Set-TextValue $ws.Range("D43") "5.68"
$ws.Range("E43").Value = "  -4.71%  "

# Row 44 - Celestia
Set-TextValue $ws.Range("D44") "11.95"
$ws.Range("E44").Value = "  -4.50%  "

# Row 45 - MultiversX
Set-TextValue $ws.Range("D45") "61.85"
$ws.Range("E45").Value = "  -6.33%  "

# Row 46 - Algorand
Set-TextValue $ws.Range("D46") "0.195"
$ws.Range("E46").Value = "  -5.51%  "

# Row 47 - FraxShare
Set-TextValue $ws.Range("D47") "8.55"
$ws.Range("E47").Value = "  -3.74%  "

# Row 48 - Cronos
Set-TextValue $ws.Range("D48") "0.0993"
$ws.Range("E48").Value = "  -3.59%  "

# Row 49 - BinanceUSD
$ws.Range("E49").Value = "  -0.13%  "

# Row 50 - ARBITRUM
Set-TextValue $ws.Range("D50") "1.16"
$ws.Range("E50").Value = "  -1.48%  "

# Row 51 - HuobiToken
$ws.Range("E51").Value = "  +5.15%  "
